# Updates cryptos list (coinranking.com scrape) - refreshed Price/Volume(1h) figures,
# and for rows 26/27 and 29/30 the two coins swapped rank order (so both the
# text/link columns and the numeric columns were replaced in place).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.147.58"
$ws.Range("E2").Value = "  -2.47%  "
# Row 3
$ws.Range("E3").Value = "  -1.57%  "
# Row 4
$ws.Range("E4").Value = "  -0.06%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.72"
$ws.Range("E5").Value = "  +0.43%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.55"
$ws.Range("E6").Value = "  -4.37%  "
# Row 7
$ws.Range("E7").Value = "  +0.07%  "
# Row 8
$ws.Range("E8").Value = "  -2.65%  "
# Row 9
$ws.Range("D9").Value = "2.435.26"
$ws.Range("E9").Value = "  -1.58%  "
# Row 10
$ws.Range("E10").Value = "  -5.29%  "
# Row 11
$ws.Range("E11").Value = "  +1.16%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.18"
$ws.Range("E12").Value = "  -2.30%  "
# Row 13
$ws.Range("E13").Value = "  -3.83%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.34"
$ws.Range("E14").Value = "  -3.42%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000171"
$ws.Range("E15").Value = "  -5.94%  "
# Row 16
$ws.Range("D16").Value = "2.868.68"
$ws.Range("E16").Value = "  -1.84%  "
# Row 17
$ws.Range("D17").Value = "62.139.41"
$ws.Range("E17").Value = "  -2.18%  "
# Row 18
$ws.Range("D18").Value = "2.432.52"
$ws.Range("E18").Value = "  -1.75%  "
# Row 19
$ws.Range("E19").Value = "  -4.96%  "
# Row 20
$ws.Range("E20").Value = "  -4.74%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "328.31"
$ws.Range("E21").Value = "  -0.97%  "
# Row 22
$ws.Range("E22").Value = "  -3.20%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.94"
$ws.Range("E23").Value = "  -8.36%  "
# Row 24
$ws.Range("E24").Value = "  -0.21%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.51"
$ws.Range("E25").Value = "  -0.79%  "
# Row 26
$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.21"
# Row 27
$ws.Range("B27").Value = "Bittensor"
$ws.Range("C27").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "627.38"
$ws.Range("E27").Value = "  -0.57%  "
# Row 28
$ws.Range("D28").Value = "2.560.32"
$ws.Range("E28").Value = "  -1.65%  "
# Row 29
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.23%  "
# Row 30
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0947"
$ws.Range("E30").Value = "  -10.27%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.42"
$ws.Range("E31").Value = "  -8.07%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.98"
$ws.Range("E32").Value = "  -5.36%  "
# Row 33
$ws.Range("E33").Value = "  -1.71%  "
# Row 34
$ws.Range("E34").Value = "  -2.21%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.92"
$ws.Range("E35").Value = "  -6.62%  "
# Row 36
$ws.Range("E36").Value = "  +0.19%  "
# Row 37
$ws.Range("E37").Value = "  -8.38%  "
# Row 38
$ws.Range("E38").Value = "  -2.59%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "149.25"
$ws.Range("E39").Value = "  +1.30%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.27"
$ws.Range("E40").Value = "  -3.50%  "
# Row 41
$ws.Range("E41").Value = "  -5.85%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.75"
$ws.Range("E42").Value = "  -4.08%  "
# Row 43
$ws.Range("E43").Value = "  +1.56%  "
# Row 44
$ws.Range("E44").Value = "  +0.00%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.47"
$ws.Range("E45").Value = "  -10.64%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "142.54"
$ws.Range("E46").Value = "  -5.68%  "
# Row 47
$ws.Range("E47").Value = "  -4.43%  "
# Row 48
$ws.Range("E48").Value = "  -3.91%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.597"
$ws.Range("E49").Value = "  -1.95%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.50"
$ws.Range("E50").Value = "  -9.76%  "
# Row 51
$ws.Range("D51").Value = "0.0₆0232"
$ws.Range("E51").Value = "  +1.62%  "
